# Add an "Address2" field to the receipt template, right after the
# existing "Address" column, shifting every following column one to the
# right (this is what Excel's "Insert Column" does when you right-click
# the header of the column that currently holds "Address" (column E) and
# choose Insert - the blank "Address" column becomes the new column and
# the old "Address" header/data slides one column to the right, into F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new, blank column at E - this shifts the old column E ("Address")
# and everything after it one column to the right.
$ws.Columns.Item(5).Insert()

# The new blank column E becomes "Address2"; the old "Address" column
# (now shifted to F) keeps its original text, so nothing else to do there.
$ws.Range("E1").Value = "Address2"

# Match the column width Excel would have copied onto the freshly inserted
# column from its left neighbour (column D).
$ws.Columns.Item(5).ColumnWidth = 9.5

# Leave the selection on the first data cell under the new header, as if
# the user had just finished typing the new header and moved down one row.
$ws.Range("F2").Select()
